$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to the new data range so that numeric-looking
# values (runs, balls, 4s, 6s, strike rate) are stored as text, matching
# the rest of the sheet (which uses numberStoredAsText).
$ws.Range("A10:K17").NumberFormat = "@"

# Row 10
$ws.Range("A10").Value = " Abu Dhabi"
$ws.Range("B10").Value = " October 30 2020"
$ws.Range("C10").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D10").Value = "Rajasthan Royals"
$ws.Range("E10").Value = "Kings XI Punjab"
$ws.Range("F10").Value = "Ben Stokes "
$ws.Range("G10").Value = "50"
$ws.Range("H10").Value = "26"
$ws.Range("I10").Value = "6"
$ws.Range("J10").Value = "3"
$ws.Range("K10").Value = "192.30"

# Row 11
$ws.Range("A11").Value = " Abu Dhabi"
$ws.Range("B11").Value = " October 25 2020"
$ws.Range("C11").Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Range("D11").Value = "Rajasthan Royals"
$ws.Range("E11").Value = "Mumbai Indians"
$ws.Range("F11").Value = "Ben Stokes "
$ws.Range("G11").Value = "107"
$ws.Range("H11").Value = "60"
$ws.Range("I11").Value = "14"
$ws.Range("J11").Value = "3"
$ws.Range("K11").Value = "178.33"

# Row 12
$ws.Range("A12").Value = " Dubai (DSC)"
$ws.Range("B12").Value = " October 14 2020"
$ws.Range("C12").Value = "Capitals won by 13 runs"
$ws.Range("D12").Value = "Rajasthan Royals"
$ws.Range("E12").Value = "Delhi Capitals"
$ws.Range("F12").Value = "Ben Stokes "
$ws.Range("G12").Value = "41"
$ws.Range("H12").Value = "35"
$ws.Range("I12").Value = "6"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "117.14"

# Row 13
$ws.Range("A13").Value = " Dubai (DSC)"
$ws.Range("B13").Value = " October 22 2020"
$ws.Range("C13").Value = "Sunrisers won by 8 wickets (with 11 balls remaining)"
$ws.Range("D13").Value = "Rajasthan Royals"
$ws.Range("E13").Value = "Sunrisers Hyderabad"
$ws.Range("F13").Value = "Ben Stokes "
$ws.Range("G13").Value = "30"
$ws.Range("H13").Value = "32"
$ws.Range("I13").Value = "2"
$ws.Range("J13").Value = "0"
$ws.Range("K13").Value = "93.75"

# Row 14
$ws.Range("A14").Value = " Dubai (DSC)"
$ws.Range("B14").Value = " November 01 2020"
$ws.Range("C14").Value = "KKR won by 60 runs"
$ws.Range("D14").Value = "Rajasthan Royals"
$ws.Range("E14").Value = "Kolkata Knight Riders"
$ws.Range("F14").Value = "Ben Stokes "
$ws.Range("G14").Value = "18"
$ws.Range("H14").Value = "11"
$ws.Range("I14").Value = "2"
$ws.Range("J14").Value = "1"
$ws.Range("K14").Value = "163.63"

# Row 15
$ws.Range("A15").Value = " Abu Dhabi"
$ws.Range("B15").Value = " October 19 2020"
$ws.Range("C15").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D15").Value = "Rajasthan Royals"
$ws.Range("E15").Value = "Chennai Super Kings"
$ws.Range("F15").Value = "Ben Stokes "
$ws.Range("G15").Value = "19"
$ws.Range("H15").Value = "11"
$ws.Range("I15").Value = "3"
$ws.Range("J15").Value = "0"
$ws.Range("K15").Value = "172.72"

# Row 16
$ws.Range("A16").Value = " Dubai (DSC)"
$ws.Range("B16").Value = " October 17 2020"
$ws.Range("C16").Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Range("D16").Value = "Rajasthan Royals"
$ws.Range("E16").Value = "Royal Challengers Bangalore"
$ws.Range("F16").Value = "Ben Stokes "
$ws.Range("G16").Value = "15"
$ws.Range("H16").Value = "19"
$ws.Range("I16").Value = "2"
$ws.Range("J16").Value = "0"
$ws.Range("K16").Value = "78.94"

# Row 17
$ws.Range("A17").Value = " Dubai (DSC)"
$ws.Range("B17").Value = " October 11 2020"
$ws.Range("C17").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D17").Value = "Rajasthan Royals"
$ws.Range("E17").Value = "Sunrisers Hyderabad"
$ws.Range("F17").Value = "Ben Stokes "
$ws.Range("G17").Value = "5"
$ws.Range("H17").Value = "6"
$ws.Range("I17").Value = "1"
$ws.Range("J17").Value = "0"
$ws.Range("K17").Value = "83.33"

